$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.611.99'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '2.304.95'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.59'
$ws.Range("E5").Value = '  +1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.75'
$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("E7").Value = '  +1.70%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.85'
$ws.Range("E10").Value = '  -0.77%  '

$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.36'
$ws.Range("E12").Value = '  +1.18%  '

$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.972'
$ws.Range("E14").Value = '  -0.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.35'
$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").Value = '2.652.45'
$ws.Range("E16").Value = '  +0.11%  '

$ws.Range("D17").Value = '2.307.21'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '42.479.59'
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.51'
$ws.Range("E19").Value = '  -2.33%  '

$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.12'
$ws.Range("E21").Value = '  -1.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.58'
$ws.Range("E22").Value = '  +2.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.05'
$ws.Range("E23").Value = '  +7.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.17'
$ws.Range("E24").Value = '  +20.30%  '

$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("E26").Value = '  -0.20%  '

$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.80'
$ws.Range("E29").Value = '  +0.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.99'
$ws.Range("E30").Value = '  +1.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.12'
$ws.Range("E31").Value = '  +1.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0877'
$ws.Range("E32").Value = '  -1.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.89'
$ws.Range("E33").Value = '  +0.99%  '

$ws.Range("E34").Value = '  +5.80%  '

$ws.Range("E35").Value = '  +1.31%  '

$ws.Range("E36").Value = '  -10.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0375'
$ws.Range("E37").Value = '  +6.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.60'
$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("E39").Value = '  +3.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.77'
$ws.Range("E40").Value = '  +1.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.50'
$ws.Range("E41").Value = '  +2.88%  '

$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.76'
$ws.Range("E42").Value = '  -1.48%  '

$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.95'
$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("E44").Value = '  -1.13%  '

$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.12'
$ws.Range("E46").Value = '  +0.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '81.16'
$ws.Range("E47").Value = '  +9.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.06'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.25'
$ws.Range("E50").Value = '  -2.32%  '

$ws.Range("D51").Value = '1.599.74'
$ws.Range("E51").Value = '  +3.88%  '
